$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 136
$ws.Range("F6").Value = 10431
$ws.Range("F8").Value = 3643
$ws.Range("F10").Value = 2481
$ws.Range("F11").Value = 51
$ws.Range("F12").Value = 2922
$ws.Range("F14").Value = 518
$ws.Range("F15").Value = 2239
$ws.Range("F17").Value = 106
$ws.Range("F19").Value = 410
$ws.Range("F22").Value = 332
$ws.Range("F23").Value = 284
$ws.Range("F24").Value = 277
$ws.Range("F25").Value = 628
$ws.Range("F26").Value = 1358
$ws.Range("F27").Value = 30
$ws.Range("F30").Value = 138
$ws.Range("F32").Value = 4012
$ws.Range("F33").Value = 3515
$ws.Range("F34").Value = 52
$ws.Range("F36").Value = 1077
$ws.Range("F37").Value = 424
$ws.Range("F40").Value = 127
$ws.Range("F41").Value = 116
$ws.Range("F42").Value = 81
$ws.Range("F45").Value = 30
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 184
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1007
$ws.Range("F5").Value = 2166
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1007
$ws.Range("F9").Value = 136
$ws.Range("F10").Value = 10431
$ws.Range("F12").Value = 3643
$ws.Range("F14").Value = 2481
$ws.Range("F15").Value = 51
$ws.Range("F16").Value = 2922
$ws.Range("F17").Value = 518
$ws.Range("F18").Value = 2239
$ws.Range("F20").Value = 106
$ws.Range("F22").Value = 410
$ws.Range("F24").Value = 332
$ws.Range("F25").Value = 277
$ws.Range("F26").Value = 628
$ws.Range("F27").Value = 1358
$ws.Range("F28").Value = 30
$ws.Range("F30").Value = 138
$ws.Range("F33").Value = 4012
$ws.Range("F34").Value = 3515
$ws.Range("F35").Value = 52
$ws.Range("F36").Value = 1077
$ws.Range("F38").Value = 424
$ws.Range("F44").Value = 127
$ws.Range("F45").Value = 81
$ws.Range("F48").Value = 30
$ws.Range("F49").Value = 184
